# Re-generate the "Chart Report" counts after re-running the combined
# report script (see commit message: "Create shell script to run all 4
# reports at once and re-generate reports").
#
# The file-type breakdown (A2:B25) is refreshed with updated counts and
# re-sorted descending by count; the four single-value summary cells
# (B28, B31, B34, B38) pick up their refreshed totals as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New file-type -> count table, already in the final (descending-by-count)
# order that the regenerated report produced.
$fileTypeRows = @(
    @("xsd",      270),
    @("docx",     142),
    @("folders",   33),
    @("html",      26),
    @("rb",        11),
    @("txt",       10),
    @("js",        10),
    @("png",        5),
    @("sh",         5),
    @("css",        5),
    @("ttf",        4),
    @("md",         2),
    @("eot",        2),
    @("svg",        2),
    @("woff",       2),
    @("woff2",      2),
    @("exe",        1),
    @("Gemfile",    1),
    @("lock",       1),
    @("bat",        1),
    @("csv",        1),
    @("xlsx",       1),
    @("otf",        1),
    @("xml",        1)
)

$startRow = 2
for ($i = 0; $i -lt $fileTypeRows.Length; $i++) {
    $row = $startRow + $i
    $label = $fileTypeRows[$i][0]
    $count = $fileTypeRows[$i][1]
    $ws.Cells.Item($row, 1).Value = $label
    $ws.Cells.Item($row, 2).Value = $count
}

# Deployed Version / Original Version counts (row 28 & 31): 125 -> 135
$ws.Cells.Item(28, 2).Value = 135
$ws.Cells.Item(31, 2).Value = 135

# Deployed Revision / Original Revision counts (row 34 & 38): 118 -> 128
$ws.Cells.Item(34, 2).Value = 128
$ws.Cells.Item(38, 2).Value = 128
